# edit.ps1
# Applies the "New crime data collected" update to the weekly CompStat (23rd
# Precinct) report: refreshes the report header (volume number + week-of dates)
# and rewrites the crime-statistics table (rows 14-29) with the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: bump volume/number and shift the reporting week by one week ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Plain numeric updates (value changes only, cell stays numeric) ---
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -66.666666666666
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 22
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = 55.555555555555
$ws.Range("L16").Value = 27.272727272727
$ws.Range("M16").Value = 7.692307692307
$ws.Range("N16").Value = -65.853658536585
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 28.571428571428
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 24
$ws.Range("I17").Value = 29
$ws.Range("J17").Value = 19
$ws.Range("K17").Value = 52.631578947368
$ws.Range("L17").Value = 93.333333333333
$ws.Range("M17").Value = 222.222222222222
$ws.Range("N17").Value = -21.621621621621
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -50
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = -60
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -80
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 3.225806451612
$ws.Range("I19").Value = 27
$ws.Range("J19").Value = 25
$ws.Range("K19").Value = 8
$ws.Range("L19").Value = 35
$ws.Range("M19").Value = 68.75
$ws.Range("N19").Value = -35.714285714285
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -37.5
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = -71.428571428571
$ws.Range("L20").Value = -71.428571428571
$ws.Range("N20").Value = -92
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -14.814814814814
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = 6.521739130434
$ws.Range("I21").Value = 77
$ws.Range("J21").Value = 73
$ws.Range("K21").Value = 5.479452054794
$ws.Range("L21").Value = 26.229508196721
$ws.Range("M21").Value = 67.391304347826
$ws.Range("N21").Value = -55.232558139534
$ws.Range("F22").Value = 1
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 30
$ws.Range("H23").Value = 16.666666666666
$ws.Range("I23").Value = 28
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = 21.739130434782
$ws.Range("L23").Value = 40
$ws.Range("M23").Value = 180
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 45.454545454545
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 50
$ws.Range("H24").Value = 10
$ws.Range("I24").Value = 36
$ws.Range("J24").Value = 40
$ws.Range("K24").Value = -10
$ws.Range("L24").Value = 16.129032258064
$ws.Range("M24").Value = -10
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -55.555555555555
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 21
$ws.Range("J25").Value = 24
$ws.Range("K25").Value = -12.5
$ws.Range("L25").Value = -25
$ws.Range("M25").Value = -44.736842105263
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = -75
$ws.Range("L26").Value = -50
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = -33.333333333333
$ws.Range("G28").Value = 2
$ws.Range("G29").Value = 2

# --- Numeric -> text transitions ("0" / "***.*" placeholders used when a ---
# --- week-over-week or 28-day comparison base is zero / not applicable)   ---
# NumberFormat is forced to Text so the digit-only placeholder "0" is not
# re-interpreted as a number; formatting is then reconciled against a cell
# that already carries the correct "text" style used throughout this table.
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E29").PasteSpecial(-4122)

# --- Text -> numeric transitions (placeholder replaced by a real figure) ---
# Formatting is reconciled against a cell that already carries the correct
# numeric style for that column (plain count vs. percentage).
$ws.Range("D20").Value = 4
$ws.Range("G14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = -100
$ws.Range("H14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = 100
$ws.Range("H14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$excel.CutCopyMode = 0
